# diodes.xlsx edit: add "DIO-GEN-0001" TPN entry on "general" sheet,
# shorten LED TPN numbering to 4 digits, fix active sheet/selection.

$wb = $excel.ActiveWorkbook
$wsGeneral = $wb.Worksheets.Item("general")
$wsLed = $wb.Worksheets.Item("LED")

# "general" sheet: A2 was a formula ("DIO-"&TEXT(ROW()-1,"000000")) producing
# "DIO-000001"; replace with the literal corrected value "DIO-GEN-0001".
$wsGeneral.Range("A2").Value = "DIO-GEN-0001"

# "LED" sheet: shrink the zero-padding in the TPN formula from 6 digits to 4
# digits (LED-000001 -> LED-0001, etc.) for rows 2-4.
$wsLed.Range("A2").Formula = '="LED-"&TEXT(ROW()-1,"0000")'
$wsLed.Range("A3").Formula = '="LED-"&TEXT(ROW()-1,"0000")'
$wsLed.Range("A4").Formula = '="LED-"&TEXT(ROW()-1,"0000")'

# Selections: "general" tab moves its selection to A3; "LED" tab selects
# A2:A4 (active cell A2) and becomes the active/selected sheet.
[void]$wsGeneral.Range("A3").Select()
[void]$wsLed.Range("A2:A4").Select()
[void]$wsLed.Activate()
